# selfTest.xlsx edit: "handling of empty param value added."
# Adds a new U15 / request-response sample row (row 15) to Sheet1, and turns
# on iterative calculation with a small max-change tolerance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# workbook.xml: <calcPr calcId="0" iterateDelta="1E-4"/>
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.0001

# New row 15: label / request JSON / empty-param response JSON
$ws.Range("A15").Value = "U15"
$ws.Range("B15").Value = '{"code":"200","data":{"total":3,"start":0,"count":3,"data":[{"supc":"SDL002940000","sellerCode":"8f7f72","name":"NM_Mobile1","pogId":621051679590,"sellerName":null,"price":0,"inventory":0,"offerPrice":0,"live":false,"imgs":null,"pageUrl":"product/nmmobile1/621051679590","category":null,"nodePath":null,"sdplus":false,"shippingDays":0,"productDesc":null,"rating":0.0,"noOfRating":0,"discount":0,"soldOut":true,"brand":"Apple","adCreated":false,"attributes":[]},{"supc":"SDL006331689","sellerCode":"8f7f72","name":"NM_Mobile1","pogId":621051679590,"sellerName":"Spice Retail Ltd.","price":2000,"inventory":200,"offerPrice":2000,"live":true,"imgs":["http://release.sdlcdn.com/http://sdstg.s3.amazonaws.com/imgs/a/a/q/NM-Mobile1-SDL006331689-1-63800.jpg"],"pageUrl":"product/nmmobile1/621051679590","category":null,"nodePath":null,"sdplus":false,"shippingDays":0,"productDesc":null,"rating":0.0,"noOfRating":0,"discount":0,"soldOut":false,"brand":"Apple","adCreated":false,"attributes":[{"name":"Color","value":"Red"}]},{"supc":"SDL000363340","sellerCode":"8f7f72","name":"NM_Mobile1","pogId":621051679590,"sellerName":"Spice Retail Ltd.","price":2000,"inventory":200,"offerPrice":2000,"live":false,"imgs":null,"pageUrl":"product/nmmobile1/621051679590","category":null,"nodePath":null,"sdplus":false,"shippingDays":0,"productDesc":null,"rating":0.0,"noOfRating":0,"discount":0,"soldOut":false,"brand":"Apple","adCreated":false,"attributes":[]}]},"status":"OK"}'
$ws.Range("C15").Value = "{}"

# C15 picks up the small 8pt black Calibri wrap-text style used for the
# "empty" response column.
$ws.Range("C15").WrapText = $true
$ws.Range("C15").Font.Size = 8
$ws.Range("C15").Font.Color = 0

# Match the author's final selection / scroll position (row 15).
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C15").Select()
